$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 15, shifting existing rows 15-19 down to 16-20
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly record
$ws.Cells.Item(15, 1).Value = 7
$ws.Cells.Item(15, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(15, 3).Value = "Ñuble"
$ws.Cells.Item(15, 4).Value = 44504
$ws.Cells.Item(15, 5).Value = 16
$ws.Cells.Item(15, 6).Value = 100112026
$ws.Cells.Item(15, 7).Value = "Haba"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 60
$ws.Cells.Item(15, 11).Value = 8000
$ws.Cells.Item(15, 12).Value = 9000
$ws.Cells.Item(15, 13).Value = 8500
$ws.Cells.Item(15, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(15, 15).Value = "Región del Maule"
$ws.Cells.Item(15, 16).Value = 340
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"

$wb.Save()
